$d = $word.ActiveDocument

$replacements = @(
    @("466×8=3728", "542×7=3794"),
    @("609×3=1827", "695×4=2780"),
    @("411×5=2055", "681×3=2043"),
    @("550×8=4400", "341×7=2387"),
    @("959×9=8631", "229×3=687"),
    @("609×6=3654", "742×9=6678"),
    @("945×3=2835", "826×4=3304"),
    @("907×4=3628", "739×3=2217"),
    @("143×8=1144", "464×4=1856"),
    @("838×3=2514", "878×2=1756"),
    @("563×5=2815", "158×7=1106"),
    @("195×6=1170", "506×2=1012"),
    @("943×7=6601", "976×5=4880"),
    @("383×3=1149", "401×6=2406"),
    @("524×5=2620", "399×4=1596"),
    @("558×7=3906", "375×3=1125"),
    @("396×4=1584", "822×7=5754"),
    @("396×5=1980", "253×7=1771"),
    @("273×8=2184", "914×6=5484"),
    @("660×2=1320", "362×5=1810"),
    @("241×5=1205", "342×4=1368"),
    @("424×2=848",  "428×4=1712"),
    @("136×6=816",  "217×8=1736"),
    @("970×6=5820", "621×5=3105"),
    @("483×3=1449", "538×4=2152")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
